# Generate Report for Handback
# Refresh the handoff/handback timestamps that get stamped when the
# handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 6e05726c-fa6f-483e-ab51-c3d5241d7e74.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 23:00:25"

# zh-cn sheet: Handoff / Handback datetimes for the same source file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-13 23:00:17"
$wsZhCn.Range("K3").Value = "2016-08-13 23:00:48"

# de-de sheet: Handback datetime for the same source file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-13 23:00:57"
